$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 15153553
$ws.Range("I15").Value = 15153553
$ws.Range("K15").Value = 45460659
$ws.Range("M15").Value = -45460490
$ws.Range("H63").Value = 24509.2
$ws.Range("I63").Value = 10000
$ws.Range("J63").Value = 28136.5
$ws.Range("K63").Value = 10000
$ws.Range("L63").Value = 28136.5
$ws.Range("M63").Value = -9376
$ws.Range("N63").Value = -29384.5
$ws.Range("H66").Value = 24509.2
$ws.Range("I66").Value = 10000
$ws.Range("J66").Value = 28136.5
$ws.Range("K66").Value = 30000
$ws.Range("L66").Value = 84409.5
$ws.Range("M66").Value = -26880
$ws.Range("N66").Value = -90649.5
$ws.Range("H107").Value = 709.2857
$ws.Range("I107").Value = 753.65216
$ws.Range("J107").Value = 655.5789
$ws.Range("K107").Value = 753.65216
$ws.Range("L107").Value = 655.5789
$ws.Range("M107").Value = 1166.34784
$ws.Range("N107").Value = -4495.5789
$ws.Range("H111").Value = 592.5789
$ws.Range("I111").Value = 450.5
$ws.Range("K111").Value = 1351.5
$ws.Range("M111").Value = 1715.5
$ws.Range("H112").Value = 14766295
$ws.Range("I112").Value = 2087.5
$ws.Range("K112").Value = 6262.5
$ws.Range("M112").Value = -5154.5
$ws.Range("H113").Value = 4693
$ws.Range("I113").Value = 2955.8
$ws.Range("J113").Value = 5658.1113
$ws.Range("K113").Value = 2955.8
$ws.Range("L113").Value = 5658.1113
$ws.Range("M113").Value = 298.1999999999998
$ws.Range("N113").Value = -12166.1113
$ws.Range("H129").Value = 4902995.5
$ws.Range("I129").Value = 50001620
$ws.Range("J129").Value = 971.3043
$ws.Range("K129").Value = 150004860
$ws.Range("L129").Value = 2913.9129
$ws.Range("M129").Value = -149999860
$ws.Range("N129").Value = -12913.9129
$ws.Range("H138").Value = 3973.8
$ws.Range("I138").Value = 1701.1316
$ws.Range("J138").Value = 5811.2764
$ws.Range("K138").Value = 5103.3948
$ws.Range("L138").Value = 17433.8292
$ws.Range("M138").Value = 36.60519999999997
$ws.Range("N138").Value = -27713.8292

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7578038.5
$ws.Range("I2").Value = 15626873
$ws.Range("J2").Value = 2664.7058
$ws.Range("K2").Value = 15626873
$ws.Range("L2").Value = 2664.7058
$ws.Range("M2").Value = -15626760
$ws.Range("N2").Value = -2890.7058
$ws.Range("H8").Value = 10001.2
$ws.Range("I8").Value = 5000
$ws.Range("J8").Value = 11251.5
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 11251.5
$ws.Range("M8").Value = -4856
$ws.Range("N8").Value = -11539.5
$ws.Range("H23").Value = 51185.273
$ws.Range("J23").Value = 40377.5
$ws.Range("L23").Value = 40377.5
$ws.Range("N23").Value = -40895.5
$ws.Range("H34").Value = 49583.43
$ws.Range("I34").Value = 7000
$ws.Range("J34").Value = 56680.668
$ws.Range("K34").Value = 7000
$ws.Range("L34").Value = 56680.668
$ws.Range("M34").Value = -6729
$ws.Range("N34").Value = -57222.668
$ws.Range("H37").Value = 12144
$ws.Range("I37").Value = 8000
$ws.Range("J37").Value = 12736
$ws.Range("K37").Value = 8000
$ws.Range("L37").Value = 12736
$ws.Range("M37").Value = -7727
$ws.Range("N37").Value = -13282
$ws.Range("H39").Value = 4861
$ws.Range("I39").Value = 4861
$ws.Range("K39").Value = 4861
$ws.Range("M39").Value = -4341
$ws.Range("H44").Value = 13029.4
$ws.Range("J44").Value = 13029.4
$ws.Range("L44").Value = 13029.4
$ws.Range("N44").Value = -14005.4
$ws.Range("H55").Value = 6051.3335
$ws.Range("J55").Value = 8053
$ws.Range("L55").Value = 8053
$ws.Range("N55").Value = -8683
$ws.Range("H80").Value = 21106.4
$ws.Range("J80").Value = 21106.4
$ws.Range("L80").Value = 21106.4
$ws.Range("N80").Value = -23102.4
$ws.Range("H83").Value = 21106.4
$ws.Range("J83").Value = 21106.4
$ws.Range("L83").Value = 63319.2
$ws.Range("N83").Value = -73303.20000000001
$ws.Range("H97").Value = 409.89285
$ws.Range("I97").Value = 367.5
$ws.Range("K97").Value = 367.5
$ws.Range("M97").Value = 128.5
$ws.Range("H110").Value = 1250.7179
$ws.Range("I110").Value = 517.7931
$ws.Range("J110").Value = 3376.2
$ws.Range("K110").Value = 517.7931
$ws.Range("L110").Value = 3376.2
$ws.Range("M110").Value = 1527.2069
$ws.Range("N110").Value = -7466.2
$ws.Range("H116").Value = 7578038.5
$ws.Range("I116").Value = 15626873
$ws.Range("J116").Value = 2664.7058
$ws.Range("K116").Value = 15626873
$ws.Range("L116").Value = 2664.7058
$ws.Range("M116").Value = -15624579
$ws.Range("N116").Value = -7252.7058

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7578038.5
$ws.Range("I3").Value = 15626873
$ws.Range("J3").Value = 2664.7058
$ws.Range("K3").Value = 15626873
$ws.Range("L3").Value = 2664.7058
$ws.Range("M3").Value = -15626759
$ws.Range("N3").Value = -2892.7058
$ws.Range("H97").Value = 13965.6
$ws.Range("I97").Value = 4914
$ws.Range("K97").Value = 4914
$ws.Range("M97").Value = -3923
$ws.Range("H105").Value = 1530.8334
$ws.Range("I105").Value = 1387
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 1387
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = 360
$ws.Range("N105").Value = -5744
$ws.Range("H107").Value = 2425
$ws.Range("I107").Value = 1668.8572
$ws.Range("J107").Value = 3748.25
$ws.Range("K107").Value = 1668.8572
$ws.Range("L107").Value = 3748.25
$ws.Range("M107").Value = 251.1428000000001
$ws.Range("N107").Value = -7588.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 25562.45
$ws.Range("I35").Value = 740.6
$ws.Range("J35").Value = 100028
$ws.Range("K35").Value = 740.6
$ws.Range("L35").Value = 100028
$ws.Range("M35").Value = -446.6
$ws.Range("N35").Value = -100616
$ws.Range("H107").Value = 1287.9333
$ws.Range("I107").Value = 702.625
$ws.Range("K107").Value = 702.625
$ws.Range("M107").Value = 1217.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 80031
$ws.Range("J44").Value = 80031
$ws.Range("L44").Value = 80031
$ws.Range("N44").Value = -81223
$ws.Range("H132").Value = 3408.3555
$ws.Range("I132").Value = 3136.8125
$ws.Range("J132").Value = 4076.7693
$ws.Range("K132").Value = 9410.4375
$ws.Range("L132").Value = 12230.3079
$ws.Range("M132").Value = -6880.4375
$ws.Range("N132").Value = -17290.3079
$ws.Range("H136").Value = 16058
$ws.Range("J136").Value = 16058
$ws.Range("L136").Value = 48174
$ws.Range("N136").Value = -53274

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 19381.4
$ws.Range("I32").Value = 1723
$ws.Range("J32").Value = 90015
$ws.Range("K32").Value = 1723
$ws.Range("L32").Value = 90015
$ws.Range("M32").Value = -1406
$ws.Range("N32").Value = -90649
$ws.Range("H40").Value = 4887.5713
$ws.Range("I40").Value = 5404.6665
$ws.Range("K40").Value = 5404.6665
$ws.Range("M40").Value = -5268.6665
$ws.Range("H136").Value = 3337617
$ws.Range("I136").Value = 5003385.5
$ws.Range("J136").Value = 6080.5
$ws.Range("K136").Value = 15010156.5
$ws.Range("L136").Value = 18241.5
$ws.Range("M136").Value = -15007606.5
$ws.Range("N136").Value = -23341.5
